$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "Ориентация на результат vs процесс"
$ws.Range("A4").Value = "Целеустремлённость и настойчивость"
